$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 278; this shifts the existing row 278 (and
# everything below it) down to row 279, growing the sheet to A1:R339.
$ws.Rows.Item(278).Insert()

# Populate the newly inserted row 278. Its content duplicates what is now
# row 279 (the old row 278), except for the Fecha (D) and Volumen (J)
# values, which take the new figures from the commit.
$ws.Cells.Item(278, 1).Value = 8
$ws.Cells.Item(278, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(278, 3).Value = "Coquimbo"
$ws.Cells.Item(278, 4).Value = 44889
$ws.Cells.Item(278, 5).Value = 4
$ws.Cells.Item(278, 6).Value = 100112012
$ws.Cells.Item(278, 7).Value = "Espinaca"
$ws.Cells.Item(278, 8).Value = "Sin especificar"
$ws.Cells.Item(278, 9).Value = "Primera"
$ws.Cells.Item(278, 10).Value = 2000
$ws.Cells.Item(278, 11).Value = 500
$ws.Cells.Item(278, 12).Value = 600
$ws.Cells.Item(278, 13).Value = 550
$ws.Cells.Item(278, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(278, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(278, 16).Value = 1100
$ws.Cells.Item(278, 17).Value = 0.5
$ws.Cells.Item(278, 18).Value = "Hortaliza"

# Match the date display format used by the rest of column D.
$ws.Cells.Item(278, 4).NumberFormat = $ws.Cells.Item(279, 4).NumberFormat
